# Zdetl Correlative Words.xlsx - add "Pronouns etc." worksheet
# -----------------------------------------------------------------
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Add a brand-new worksheet at the end of the workbook (after
#    "Verb conjugations and endings") and name it.
# ---------------------------------------------------------------
$sheetCount = $wb.Worksheets.Count
$lastSheet  = $wb.Worksheets.Item($sheetCount)
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws3.Name = "Pronouns etc."

# ---------------------------------------------------------------
# 2. Title (row 1)
# ---------------------------------------------------------------
$ws3.Range("A1").Value = "Pronouns, Articles, and other parts of speech"
$ws3.Range("A1").Font.Bold = $true
$ws3.Range("A1").Font.Size = 16
$ws3.Rows.Item(1).RowHeight = 21

# ---------------------------------------------------------------
# 3. "Standard forms" pronoun table (rows 2-9)
# ---------------------------------------------------------------
$ws3.Range("A2").Value = "Standard forms"
$ws3.Range("E2").Value = "Possessive"
$ws3.Range("F2").Value = "add '-o'"

$ws3.Range("A3").Value = "I"
$ws3.Range("B3").Value = "first person singular"
$ws3.Range("C3").Value = "ze"
$ws3.Range("E3").Value = "my"
$ws3.Range("F3").Value = "zeo"

$ws3.Range("A4").Value = "we"
$ws3.Range("B4").Value = "first person plural"
$ws3.Range("C4").Value = "de"
$ws3.Range("E4").Value = "our"
$ws3.Range("F4").Value = "deo"

$ws3.Range("A5").Value = "you"
$ws3.Range("B5").Value = "second person singular"
$ws3.Range("C5").Value = "ve"
$ws3.Range("E5").Value = "your"
$ws3.Range("F5").Value = "veo"

$ws3.Range("A6").Value = "y'all"
$ws3.Range("B6").Value = "second person plural"
$ws3.Range("C6").Value = "vi"
$ws3.Range("E6").Value = "y'all's"
$ws3.Range("F6").Value = "vio"

$ws3.Range("A7").Value = "he/she/it"
$ws3.Range("B7").Value = "third person singular"
$ws3.Range("C7").Value = "se"
$ws3.Range("E7").Value = "his/hers/its"
$ws3.Range("F7").Value = "seo"
$ws3.Range("G7").Value = 'Note: this can also be interpreted as the singular "they" in Anglic.'

$ws3.Range("A8").Value = "they"
$ws3.Range("B8").Value = "third person plural"
$ws3.Range("C8").Value = "ye"
$ws3.Range("E8").Value = "theirs"
$ws3.Range("F8").Value = "yeo"

$ws3.Range("A9").Value = "reflexive pronoun"
$ws3.Range("B9").Value = 'non-specific "one"'
$ws3.Range("C9").Value = "si"
$ws3.Range("E9").Value = "one's"
$ws3.Range("F9").Value = "sio"

# ---------------------------------------------------------------
# 4. Definite article (rows 11-12)
# ---------------------------------------------------------------
$ws3.Range("A11").Value = "definite article"
$ws3.Range("B11").Value = "the"
$ws3.Range("C11").Value = "ke"
$ws3.Range("A12").Value = "there is no indefinite article in Zdetl."

# ---------------------------------------------------------------
# 5. Misc notes (rows 14-16)
# ---------------------------------------------------------------
$ws3.Range("A14").Value = 'Verb infinitives typically end in "-e^"'
$ws3.Range("A15").Value = 'Adverbs typically end in "-ie"'
$ws3.Range("A16").Value = 'Adjectives typically end in "-o"'

# ---------------------------------------------------------------
# 6. Final note (row 18)
# ---------------------------------------------------------------
$ws3.Range("A18").Value = "Nouns don't always follow a standard pattern"

# ---------------------------------------------------------------
# 7. Column widths (approximate best-fit)
# ---------------------------------------------------------------
$ws3.Columns.Item(1).ColumnWidth = 15.44140625
$ws3.Columns.Item(2).ColumnWidth = 19.6640625
$ws3.Columns.Item(5).ColumnWidth = 10.109375

# ---------------------------------------------------------------
# 8. Page setup & final selection / view state
# ---------------------------------------------------------------
$ws3.PageSetup.Orientation = 1
$ws3.Range("A18").Select() | Out-Null

Write-Output "Added 'Pronouns etc.' worksheet"
